$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated LR-pair stats (Cxcl12-Cxcr4), recomputed "following Dr Hou advice":
# existing rows 2-13 get new values, and 4 new rows (14-17, sending cluster "sCs")
# are appended, each row = [A..T] = 20 columns.
$data = @(
    @("ECs", "Cxcl12", "Cxcr4", "ECs", 3, 1, 81.05837566666666, 243.175127, 0.3545816884225585, 0.3545816884225585, 3, 1, 161.9042856666667, 485.712857, 0.3829374249648381, 0.3829374249648381, 13123.69840961198, 118113.2856865078, 0.1357825987042191, 0.1357825987042191),
    @("ECs", "Cxcl12", "Cxcr4", "FAPs", 3, 1, 81.05837566666666, 243.175127, 0.3545816884225585, 0.3545816884225585, 2, 0.6666666666666666, 0.307744, 0.9232319999999999, 0.0007278787860563825, 0.0007278787860563825, 24.94522876116266, 224.5070588504639, 0.0002580924889268343, 0.0002580924889268343),
    @("ECs", "Cxcl12", "Cxcr4", "M2", 3, 1, 81.05837566666666, 243.175127, 0.3545816884225585, 0.3545816884225585, 3, 1, 225.0171966666667, 675.05159, 0.5322126311204886, 0.5322126311204886, 18239.52845886688, 164155.7561298019, 0.1887128533425151, 0.1887128533425151),
    @("ECs", "Cxcl12", "Cxcr4", "sCs", 3, 1, 81.05837566666666, 243.175127, 0.3545816884225585, 0.3545816884225585, 3, 1, 35.56644499999999, 106.699335, 0.08412206512861695, 0.08412206512861695, 2882.958259937838, 25946.62433944054, 0.02982814388689743, 0.02982814388689743),
    @("FAPs", "Cxcl12", "Cxcr4", "ECs", 3, 1, 123.018252, 369.054756, 0.5381309351710768, 0.5381309351710768, 3, 1, 161.9042856666667, 485.712857, 0.3829374249648381, 0.3829374249648381, 19917.18221402199, 179254.6399261979, 0.2060704746083324, 0.2060704746083324),
    @("FAPs", "Cxcl12", "Cxcr4", "FAPs", 3, 1, 123.018252, 369.054756, 0.5381309351710768, 0.5381309351710768, 2, 0.6666666666666666, 0.307744, 0.9232319999999999, 0.0007278787860563825, 0.0007278787860563825, 37.858128943488, 340.723160491392, 0.0003916940918317092, 0.0003916940918317092),
    @("FAPs", "Cxcl12", "Cxcr4", "M2", 3, 1, 123.018252, 369.054756, 0.5381309351710768, 0.5381309351710768, 3, 1, 225.0171966666667, 675.05159, 0.5322126311204886, 0.5322126311204886, 27681.22220387356, 249130.999834862, 0.2864000808947278, 0.2864000808947278),
    @("FAPs", "Cxcl12", "Cxcr4", "sCs", 3, 1, 123.018252, 369.054756, 0.5381309351710768, 0.5381309351710768, 3, 1, 35.56644499999999, 106.699335, 0.08412206512861695, 0.08412206512861695, 4375.32189375414, 39377.89704378726, 0.04526868557618487, 0.04526868557618487),
    @("M2", "Cxcl12", "Cxcr4", "ECs", 3, 1, 0.3624666666666667, 1.0874, 0.001585573873230423, 0.001585573873230423, 3, 1, 161.9042856666667, 485.712857, 0.3829374249648381, 0.3829374249648381, 58.68490674464445, 528.1641607018, 0.0006071755761063829, 0.0006071755761063829),
    @("M2", "Cxcl12", "Cxcr4", "FAPs", 3, 1, 0.3624666666666667, 1.0874, 0.001585573873230423, 0.001585573873230423, 2, 0.6666666666666666, 0.307744, 0.9232319999999999, 0.0007278787860563825, 0.0007278787860563825, 0.1115469418666667, 1.0039224768, 0.000001154105586049677, 0.000001154105586049677),
    @("M2", "Cxcl12", "Cxcr4", "M2", 3, 1, 0.3624666666666667, 1.0874, 0.001585573873230423, 0.001585573873230423, 3, 1, 225.0171966666667, 675.05159, 0.5322126311204886, 0.5322126311204886, 81.56123321844444, 734.0510989659999, 0.0008438624429078676, 0.0008438624429078676),
    @("M2", "Cxcl12", "Cxcr4", "sCs", 3, 1, 0.3624666666666667, 1.0874, 0.001585573873230423, 0.001585573873230423, 3, 1, 35.56644499999999, 106.699335, 0.08412206512861695, 0.08412206512861695, 12.89165076433333, 116.024856879, 0.0001333817486301231, 0.0001333817486301231),
    @("sCs", "Cxcl12", "Cxcr4", "ECs", 3, 1, 24.16373066666667, 72.491192, 0.1057018025331343, 0.1057018025331344, 3, 1, 161.9042856666667, 485.712857, 0.3829374249648381, 0.3829374249648381, 3912.211552628394, 35209.90397365555, 0.04047717607618026, 0.04047717607618027),
    @("sCs", "Cxcl12", "Cxcr4", "FAPs", 3, 1, 24.16373066666667, 72.491192, 0.1057018025331343, 0.1057018025331344, 2, 0.6666666666666666, 0.307744, 0.9232319999999999, 0.0007278787860563825, 0.0007278787860563825, 7.436243130282666, 66.92618817254399, 0.00007693809971178928, 0.00007693809971178928),
    @("sCs", "Cxcl12", "Cxcr4", "M2", 3, 1, 24.16373066666667, 72.491192, 0.1057018025331343, 0.1057018025331344, 3, 1, 225.0171966666667, 675.05159, 0.5322126311204886, 0.5322126311204886, 5437.254935621698, 48935.29442059528, 0.05625583444033775, 0.05625583444033776),
    @("sCs", "Cxcl12", "Cxcr4", "sCs", 3, 1, 24.16373066666667, 72.491192, 0.1057018025331343, 0.1057018025331344, 3, 1, 35.56644499999999, 106.699335, 0.08412206512861695, 0.08412206512861695, 859.4179977508131, 7734.761979757319, 0.008891853916904535, 0.008891853916904537)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $rowVals = $data[$i]
    $r = $startRow + $i
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowVals[$c]
    }
}
